$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '45.135.59'
$ws.Cells.Item(2, 5).Value = '  -1.32%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '2.357.06'
$ws.Cells.Item(3, 5).Value = '  -1.86%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.09%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '330.97'
$ws.Cells.Item(5, 5).Value = '  +3.35%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '106.96'
$ws.Cells.Item(6, 5).Value = '  -7.31%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.634'
$ws.Cells.Item(7, 5).Value = '  -0.74%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.18%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.614'
$ws.Cells.Item(9, 5).Value = '  -2.78%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '40.81'
$ws.Cells.Item(10, 5).Value = '  -4.93%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.0917'
$ws.Cells.Item(11, 5).Value = '  -1.84%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '8.42'
$ws.Cells.Item(12, 5).Value = '  -3.77%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  -1.02%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '0.973'
$ws.Cells.Item(14, 5).Value = '  -4.58%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '2.722.30'
$ws.Cells.Item(15, 5).Value = '  -1.57%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '15.32'
$ws.Cells.Item(16, 5).Value = '  -4.50%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '2.344.30'
$ws.Cells.Item(17, 5).Value = '  -2.77%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '45.171.30'
$ws.Cells.Item(18, 5).Value = '  -1.28%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '15.16'
$ws.Cells.Item(19, 5).Value = '  +11.07%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '7.27'
$ws.Cells.Item(20, 5).Value = '  -3.75%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  -2.54%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '3.68'
$ws.Cells.Item(22, 5).Value = '  +3.06%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '72.95'
$ws.Cells.Item(23, 5).Value = '  -3.10%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '258.74'
$ws.Cells.Item(24, 5).Value = '  -2.73%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '2.29'
$ws.Cells.Item(25, 5).Value = '  -3.78%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +0.03%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '11.29'
$ws.Cells.Item(27, 5).Value = '  -1.09%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '7.40'
$ws.Cells.Item(28, 5).Value = '  -3.17%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '2.29'
$ws.Cells.Item(29, 5).Value = '  -2.60%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '0.0962'
$ws.Cells.Item(30, 5).Value = '  -4.56%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '22.20'
$ws.Cells.Item(31, 5).Value = '  -3.03%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '36.71'
$ws.Cells.Item(32, 5).Value = '  -9.31%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '166.91'
$ws.Cells.Item(33, 5).Value = '  -3.65%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '2.83'
$ws.Cells.Item(34, 5).Value = '  -4.54%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '3.30'
$ws.Cells.Item(35, 5).Value = '  +5.15%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'Stellar'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.132'
$ws.Cells.Item(36, 5).Value = '  -1.02%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  -3.62%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '4.71'
$ws.Cells.Item(38, 5).Value = '  -6.62%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '1.93'
$ws.Cells.Item(39, 5).Value = '  +8.62%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '3.97'
$ws.Cells.Item(40, 5).Value = '  -6.56%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.0353'
$ws.Cells.Item(41, 5).Value = '  -3.73%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '97.11'
$ws.Cells.Item(42, 5).Value = '  -3.36%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '69.72'
$ws.Cells.Item(43, 5).Value = '  -3.94%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '1.878.57'
$ws.Cells.Item(44, 5).Value = '  +13.56%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'THORChain'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '6.08'
$ws.Cells.Item(45, 5).Value = '  +2.88%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'Algorand'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.228'
$ws.Cells.Item(46, 5).Value = '  -7.29%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '12.89'
$ws.Cells.Item(47, 5).Value = '  -6.86%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '1.00'
$ws.Cells.Item(48, 5).Value = '  +0.25%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'ordi'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '84.75'
$ws.Cells.Item(49, 5).Value = '  -5.07%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '111.74'
$ws.Cells.Item(50, 5).Value = '  -4.39%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '9.29'
$ws.Cells.Item(51, 5).Value = '  -2.25%  '
